$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write values for the new rows (358-366)
$ws.Range("A358").Value = 44432
$ws.Range("B358").Value = 4
$ws.Range("C358").Value = 11
$ws.Range("D358").Value = 118.4834123222749

$ws.Range("A359").Value = 44433
$ws.Range("B359").Value = 1
$ws.Range("C359").Value = 12
$ws.Range("D359").Value = 129.2546316242999

$ws.Range("A360").Value = 44434
$ws.Range("B360").Value = 2
$ws.Range("C360").Value = 10
$ws.Range("D360").Value = 107.7121930202499

$ws.Range("A361").Value = 44435
$ws.Range("B361").Value = 1
$ws.Range("C361").Value = 9
$ws.Range("D361").Value = 96.9409737182249

$ws.Range("A362").Value = 44436
$ws.Range("B362").Value = 2
$ws.Range("C362").Value = 11
$ws.Range("D362").Value = 118.4834123222749

$ws.Range("A363").Value = 44437
$ws.Range("B363").Value = 5
$ws.Range("C363").Value = 15
$ws.Range("D363").Value = 161.5682895303748

$ws.Range("A364").Value = 44438
$ws.Range("B364").Value = 1
$ws.Range("C364").Value = 16
$ws.Range("D364").Value = 172.3395088323998

$ws.Range("A365").Value = 44439
$ws.Range("B365").Value = 0
$ws.Range("C365").Value = 12
$ws.Range("D365").Value = 129.2546316242999

$ws.Range("A366").Value = 44440
$ws.Range("B366").Value = 1
$ws.Range("C366").Value = 12
$ws.Range("D366").Value = 129.2546316242999

# Copy formatting (incl. date style on column A) from the last existing row
$ws.Range("A357:D357").Copy()
$ws.Range("A358:D366").PasteSpecial(-4122)

$excel.CutCopyMode = 0

Write-Output "Added rows 358 to 366"
